$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-07 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-08 Saturday", 2) | Out-Null
$d.Content.Find.Execute("918÷2=459, 0", $true, $false, $false, $false, $false, $true, 1, $false, "433÷9=48, 1", 2) | Out-Null
$d.Content.Find.Execute("990÷3=330, 0", $true, $false, $false, $false, $false, $true, 1, $false, "713÷5=142, 3", 2) | Out-Null
$d.Content.Find.Execute("922÷3=307, 1", $true, $false, $false, $false, $false, $true, 1, $false, "438÷3=146, 0", 2) | Out-Null
$d.Content.Find.Execute("835÷2=417, 1", $true, $false, $false, $false, $false, $true, 1, $false, "448÷7=64, 0", 2) | Out-Null
$d.Content.Find.Execute("170÷7=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "672÷6=112, 0", 2) | Out-Null
$d.Content.Find.Execute("980÷7=140, 0", $true, $false, $false, $false, $false, $true, 1, $false, "434÷2=217, 0", 2) | Out-Null
$d.Content.Find.Execute("561÷5=112, 1", $true, $false, $false, $false, $false, $true, 1, $false, "602÷4=150, 2", 2) | Out-Null
$d.Content.Find.Execute("402÷6=67, 0", $true, $false, $false, $false, $false, $true, 1, $false, "384÷4=96, 0", 2) | Out-Null
$d.Content.Find.Execute("391÷4=97, 3", $true, $false, $false, $false, $false, $true, 1, $false, "914÷7=130, 4", 2) | Out-Null
$d.Content.Find.Execute("501÷9=55, 6", $true, $false, $false, $false, $false, $true, 1, $false, "494÷4=123, 2", 2) | Out-Null
$d.Content.Find.Execute("584÷6=97, 2", $true, $false, $false, $false, $false, $true, 1, $false, "110÷6=18, 2", 2) | Out-Null
$d.Content.Find.Execute("319÷5=63, 4", $true, $false, $false, $false, $false, $true, 1, $false, "327÷6=54, 3", 2) | Out-Null
$d.Content.Find.Execute("450÷2=225, 0", $true, $false, $false, $false, $false, $true, 1, $false, "681÷8=85, 1", 2) | Out-Null
$d.Content.Find.Execute("823÷9=91, 4", $true, $false, $false, $false, $false, $true, 1, $false, "234÷6=39, 0", 2) | Out-Null
$d.Content.Find.Execute("701÷2=350, 1", $true, $false, $false, $false, $false, $true, 1, $false, "355÷9=39, 4", 2) | Out-Null
$d.Content.Find.Execute("134÷5=26, 4", $true, $false, $false, $false, $false, $true, 1, $false, "179÷9=19, 8", 2) | Out-Null
$d.Content.Find.Execute("631÷8=78, 7", $true, $false, $false, $false, $false, $true, 1, $false, "419÷3=139, 2", 2) | Out-Null
$d.Content.Find.Execute("111÷8=13, 7", $true, $false, $false, $false, $false, $true, 1, $false, "939÷3=313, 0", 2) | Out-Null
$d.Content.Find.Execute("962÷4=240, 2", $true, $false, $false, $false, $false, $true, 1, $false, "950÷7=135, 5", 2) | Out-Null
$d.Content.Find.Execute("151÷7=21, 4", $true, $false, $false, $false, $false, $true, 1, $false, "810÷2=405, 0", 2) | Out-Null
$d.Content.Find.Execute("710÷4=177, 2", $true, $false, $false, $false, $false, $true, 1, $false, "525÷2=262, 1", 2) | Out-Null
$d.Content.Find.Execute("492÷2=246, 0", $true, $false, $false, $false, $false, $true, 1, $false, "906÷7=129, 3", 2) | Out-Null
$d.Content.Find.Execute("550÷7=78, 4", $true, $false, $false, $false, $false, $true, 1, $false, "185÷5=37, 0", 2) | Out-Null
$d.Content.Find.Execute("326÷7=46, 4", $true, $false, $false, $false, $false, $true, 1, $false, "984÷9=109, 3", 2) | Out-Null
$d.Content.Find.Execute("790÷9=87, 7", $true, $false, $false, $false, $false, $true, 1, $false, "405÷2=202, 1", 2) | Out-Null
